# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Mon Feb 19 11:52:25 UTC 2024 with GitHub Actions"
#
# Column D ("Price") values are stored as plain text in the source sheet (not
# numbers), including values that look numeric (e.g. "40.00", "0.630"). A bare
# assignment of a numeric-looking string lets Excel auto-convert the cell to a
# Number, silently dropping meaningful trailing zeros / thousands-dot formatting
# (e.g. "120.90" -> 120.9, or "52.416.55" would misparse). Prefixing with a
# leading apostrophe keeps every Price cell as literal Text, matching the source.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'52.416.55"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3
$ws.Range("D3").Value = "'2.914.78"
$ws.Range("E3").Value = "  +3.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'353.15"
$ws.Range("E5").Value = "  -0.42%  "

# Row 6
$ws.Range("D6").Value = "'112.97"
$ws.Range("E6").Value = "  +0.87%  "

# Row 7
$ws.Range("E7").Value = "  +0.84%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("D9").Value = "'0.630"
$ws.Range("E9").Value = "  -0.47%  "

# Row 10
$ws.Range("D10").Value = "'40.00"
$ws.Range("E10").Value = "  -0.92%  "

# Row 11
$ws.Range("E11").Value = "  +3.46%  "

# Row 12
$ws.Range("E12").Value = "  +0.25%  "

# Row 13
$ws.Range("D13").Value = "'19.92"
$ws.Range("E13").Value = "  -0.57%  "

# Row 14
$ws.Range("E14").Value = "  +0.64%  "

# Row 15
$ws.Range("D15").Value = "'3.371.92"
$ws.Range("E15").Value = "  +3.93%  "

# Row 16
$ws.Range("E16").Value = "  +6.34%  "

# Row 17
$ws.Range("D17").Value = "'2.899.31"
$ws.Range("E17").Value = "  +3.22%  "

# Row 18
$ws.Range("D18").Value = "'52.404.75"
$ws.Range("E18").Value = "  +1.21%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'7.64"
$ws.Range("E19").Value = "  -0.22%  "

# Row 20
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").Value = "'3.31"
$ws.Range("E20").Value = "  +3.81%  "

# Row 21
$ws.Range("D21").Value = "'14.23"
$ws.Range("E21").Value = "  +3.99%  "

# Row 22
$ws.Range("E22").Value = "  +0.48%  "

# Row 23
$ws.Range("D23").Value = "'70.97"
$ws.Range("E23").Value = "  +0.60%  "

# Row 24
$ws.Range("D24").Value = "'270.33"
$ws.Range("E24").Value = "  +0.59%  "

# Row 25
$ws.Range("E25").Value = "  +0.52%  "

# Row 26
$ws.Range("D26").Value = "'0.174"
$ws.Range("E26").Value = "  +7.52%  "

# Row 27
$ws.Range("E27").Value = "  +2.40%  "

# Row 28
$ws.Range("E28").Value = "  -0.10%  "

# Row 29
$ws.Range("E29").Value = "  +2.72%  "

# Row 30
$ws.Range("D30").Value = "'6.71"
$ws.Range("E30").Value = "  +9.32%  "

# Row 31
$ws.Range("D31").Value = "'38.17"
$ws.Range("E31").Value = "  -2.53%  "

# Row 32
$ws.Range("D32").Value = "'6.36"
$ws.Range("E32").Value = "  +12.84%  "

# Row 33
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0998"
$ws.Range("E33").Value = "  +12.67%  "

# Row 34
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.26"
$ws.Range("E34").Value = "  +0.35%  "

# Row 35
$ws.Range("D35").Value = "'53.32"
$ws.Range("E35").Value = "  +2.08%  "

# Row 36
$ws.Range("D36").Value = "'0.0453"
$ws.Range("E36").Value = "  +1.80%  "

# Row 37
$ws.Range("E37").Value = "  -0.10%  "

# Row 38
$ws.Range("E38").Value = "  +6.58%  "

# Row 39
$ws.Range("D39").Value = "'19.04"
$ws.Range("E39").Value = "  +0.52%  "

# Row 40
$ws.Range("E40").Value = "  +3.65%  "

# Row 41
$ws.Range("E41").Value = "  +13.57%  "

# Row 42
$ws.Range("E42").Value = "  +1.59%  "

# Row 43
$ws.Range("D43").Value = "'23.46"
$ws.Range("E43").Value = "  +7.05%  "

# Row 44
$ws.Range("D44").Value = "'120.90"
$ws.Range("E44").Value = "  +0.66%  "

# Row 45
$ws.Range("E45").Value = "  +7.51%  "

# Row 46
$ws.Range("E46").Value = "  -1.65%  "

# Row 47
$ws.Range("E47").Value = "  +4.18%  "

# Row 48
$ws.Range("D48").Value = "'2.197.54"
$ws.Range("E48").Value = "  +3.68%  "

# Row 49
$ws.Range("D49").Value = "'0.264"
$ws.Range("E49").Value = "  +22.35%  "

# Row 50
$ws.Range("D50").Value = "'0.0346"
$ws.Range("E50").Value = "  +12.32%  "

# Row 51
$ws.Range("D51").Value = "'0.968"
$ws.Range("E51").Value = "  +1.23%  "
